$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 is updated with a new day's electricity price data (automatic update)

$ws.Range("A2").Value = 46076

$ws.Range("B2").Value = 23.89
$ws.Range("C2").Value = 19.91
$ws.Range("D2").Value = 18.24
$ws.Range("E2").Value = 12.63
$ws.Range("F2").Value = 10.58
$ws.Range("G2").Value = 18.03
$ws.Range("H2").Value = 24.23
$ws.Range("I2").Value = 35.02
$ws.Range("J2").Value = 35.49
$ws.Range("K2").Value = 14.52
$ws.Range("L2").Value = 10.54
$ws.Range("M2").Value = 10
$ws.Range("N2").Value = 9.99
$ws.Range("O2").Value = 9.039999999999999
$ws.Range("P2").Value = 8.24
$ws.Range("Q2").Value = 10
$ws.Range("R2").Value = 13.9
$ws.Range("S2").Value = 22.03
$ws.Range("T2").Value = 37.29
$ws.Range("U2").Value = 93.01000000000001
$ws.Range("V2").Value = 175.13
$ws.Range("W2").Value = 110.76
$ws.Range("X2").Value = 54.16
$ws.Range("Y2").Value = 31.38
$ws.Range("Z2").Value = 33.67

# AA2 unchanged: "20h-24h"
$ws.Range("AB2").Value = 92.86
# AC2 unchanged: "20h-22h"
$ws.Range("AD2").Value = 142.94
$ws.Range("AE2").Value = "18h-20h"
$ws.Range("AF2").Value = 65.15000000000001
$ws.Range("AG2").Value = "0h-23h"
